$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32 becomes dogwifhat, Row 33 becomes Bittensor (swap with updated values)
$ws.Range("B32").Value = "dogwifhat"
$ws.Range("C32").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.67"
$ws.Range("E32").Value = "  +11.63%  "

$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "555.16"
$ws.Range("E33").Value = "  -2.30%  "

# Price and Volume(1h) updates
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.235.82"
$ws.Range("E2").Value = "  +2.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.386.40"
$ws.Range("E3").Value = "  +1.59%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.80"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.06"
$ws.Range("E6").Value = "  +2.59%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.596"
$ws.Range("E8").Value = "  +1.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.196"
$ws.Range("E9").Value = "  +7.63%  "
$ws.Range("E10").Value = "  +2.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.58"
$ws.Range("E11").Value = "  +3.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000282"
$ws.Range("E12").Value = "  +3.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "680.29"
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.64"
$ws.Range("E14").Value = "  +2.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.924.40"
$ws.Range("E15").Value = "  +1.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.302.15"
$ws.Range("E16").Value = "  +2.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.403.82"
$ws.Range("E17").Value = "  +1.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.73"
$ws.Range("E19").Value = "  +0.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.31"
$ws.Range("E20").Value = "  +2.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.904"
$ws.Range("E21").Value = "  +1.21%  "
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.13"
$ws.Range("E23").Value = "  +1.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "102.88"
$ws.Range("E24").Value = "  +1.23%  "
$ws.Range("E25").Value = "  +0.61%  "
$ws.Range("E26").Value = "  +1.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.61"
$ws.Range("E27").Value = "  +1.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.90"
$ws.Range("E28").Value = "  +2.83%  "
$ws.Range("E29").Value = "  +2.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.96"
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("E31").Value = "  +1.24%  "
$ws.Range("E34").Value = "  +1.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.54"
$ws.Range("E35").Value = "  +2.13%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.675.75"
$ws.Range("E37").Value = "  -1.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "35.83"
$ws.Range("E38").Value = "  +2.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.140"
$ws.Range("E39").Value = "  +4.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0714"
$ws.Range("E40").Value = "  +6.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.26"
$ws.Range("E41").Value = "  +3.47%  "
$ws.Range("E42").Value = "  +2.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.340"
$ws.Range("E43").Value = "  +1.65%  "
$ws.Range("E44").Value = "  +3.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.32"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("E46").Value = "  +1.31%  "
$ws.Range("E47").Value = "  +1.32%  "
$ws.Range("E48").Value = "  +5.69%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.25"
$ws.Range("E50").Value = "  +0.97%  "
$ws.Range("E51").Value = "  +4.29%  "
